$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 174 (Florida King / Primera row),
# pushing the existing rows 174-180 down to 176-182.
$ws.Rows.Item(174).Insert()
$ws.Rows.Item(174).Insert()

# Populate the first new row (174) with the "Doctor Davis - Primera" record.
$ws.Range("A174").Value = 11
$ws.Range("B174").Value = "Vega Monumental Concepción"
$ws.Range("C174").Value = "Bíobío"
$ws.Range("D174").Value = 44610
$ws.Range("E174").Value = 8
$ws.Range("F174").Value = "Fruta"
$ws.Range("G174").Value = 100103
$ws.Range("H174").Value = "Frutos de hueso (carozo)"
$ws.Range("I174").Value = 100103004
$ws.Range("J174").Value = "Durazno"
$ws.Range("K174").Value = "Doctor Davis"
$ws.Range("L174").Value = "Primera"
$ws.Range("M174").Value = 220
$ws.Range("N174").Value = 14000
$ws.Range("O174").Value = 15000
$ws.Range("P174").Value = 14545
$ws.Range("Q174").Value = "$/caja 16 kilos empedrada"
$ws.Range("R174").Value = "Región de O'Higgins"
$ws.Range("S174").Value = 909
$ws.Range("T174").Value = 16

# Populate the second new row (175) with the "Doctor Davis - Segunda" record.
$ws.Range("A175").Value = 11
$ws.Range("B175").Value = "Vega Monumental Concepción"
$ws.Range("C175").Value = "Bíobío"
$ws.Range("D175").Value = 44610
$ws.Range("E175").Value = 8
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100103
$ws.Range("H175").Value = "Frutos de hueso (carozo)"
$ws.Range("I175").Value = 100103004
$ws.Range("J175").Value = "Durazno"
$ws.Range("K175").Value = "Doctor Davis"
$ws.Range("L175").Value = "Segunda"
$ws.Range("M175").Value = 250
$ws.Range("N175").Value = 11000
$ws.Range("O175").Value = 12000
$ws.Range("P175").Value = 11400
$ws.Range("Q175").Value = "$/caja 16 kilos empedrada"
$ws.Range("R175").Value = "Región de O'Higgins"
$ws.Range("S175").Value = 712
$ws.Range("T175").Value = 16
